$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dates are Excel serial numbers already matching column A's
# existing values, so we keep them as plain numbers and copy the source
# row's formatting onto the new rows).
$data = @(
    @(44313, 0, 7, 46.39756081394578),
    @(44314, 1, 8, 53.02578378736661),
    @(44315, 1, 7, 46.39756081394578),
    @(44316, 1, 8, 53.02578378736661),
    @(44317, 4, 12, 79.53867568104991),
    @(44318, 1, 10, 66.28222973420826)
)

$startRow = 239
$endRow = 244

# Copy formatting from the last existing row (238) down onto the new block
# in one shot so the new cells keep the same styles (date format/border/
# alignment on column A, plain numbers on B:D) as the rest of the table.
$ws.Range("A238:D238").Copy() | Out-Null
$ws.Range("A$startRow`:D$endRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
